$wb = $excel.ActiveWorkbook
$wb.BreakLink("/Users/hecvasro/workspace/pucmm/isc-434-t/ppr/tspi/plan.xlsx")
